$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NroSiniestro value for the preproduccion row (F3):
# old: "0420172010222  " -> new: "0420172010458  " (trailing spaces preserved)
# Leading apostrophe forces Excel to keep it as text (matching original quotePrefix style)
# without altering the cell's existing style/number format.
$ws.Range("F3").Value = "'0420172010458  "

# Update the active selection shown in the sheet view from D3 to H3
$ws.Range("H3").Select()
